$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12-14: rotate the three ID values
$ws.Range("A12").Value = "103/06-29-057-03W4/00"
$ws.Range("A13").Value = "100/08-29-057-03W4/00"
$ws.Range("A14").Value = "100/07-29-057-03W4/00"

# Rows 21-22: swap the two ID values
$ws.Range("A21").Value = "100/08-21-062-23W5/00"
$ws.Range("A22").Value = "100/08-21-062-23W5/02"
